# This workbook is a weekly NFL betting "tool": the top block (rows 2-5)
# holds this week's game inputs (team, spread, total), and the two side
# lists (L2:L12 "win" scores and S2:S12 "loss" scores) are pasted-in
# reference data used by lookups further down the sheet. Updating it for
# a new week means re-entering both blocks; everything else (E:X) is
# formula-driven and recalculates automatically (including the G2:G5 /
# J2:J5 one-variable Data Tables).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: New England
$ws.Range("B2").Value = "NE"
$ws.Range("C2").Value = -4
$ws.Range("D2").Value = 34

# Row 3: Kansas City
$ws.Range("B3").Value = "KC"
$ws.Range("C3").Value = -9
$ws.Range("D3").Value = 44

# Row 4: Minnesota
$ws.Range("B4").Value = "MIN"
$ws.Range("C4").Value = -4
$ws.Range("D4").Value = 45

# Row 5: LA Rams (spread/total unchanged this week)
$ws.Range("B5").Value = "LAR"

# "W ex R" reference list (L2:L12) - updated win-score pairs
$ws.Range("L2").Value = "  27   66"
$ws.Range("L3").Value = "  20   58"
$ws.Range("L4").Value = "  24   51"
$ws.Range("L5").Value = "  31   49"
$ws.Range("L6").Value = "  30   38"
$ws.Range("L7").Value = "  23   37"
$ws.Range("L8").Value = "  34   34"
$ws.Range("L9").Value = "  17   30"
$ws.Range("L10").Value = "  28   28"
$ws.Range("L11").Value = "  19   26"
$ws.Range("L12").Value = "  38   26"

# "L ex R" reference list (S2:S12) - updated loss-score pairs
$ws.Range("S2").Value = "  17   79"
$ws.Range("S3").Value = "  10   69"
$ws.Range("S4").Value = "  16   55"
$ws.Range("S5").Value = "  13   48"
$ws.Range("S6").Value = "  20   48"
$ws.Range("S7").Value = "  14   40"
$ws.Range("S8").Value = "  24   36"
$ws.Range("S9").Value = "  3    31"
$ws.Range("S10").Value = "  9    26"
$ws.Range("S11").Value = "  21   26"
$ws.Range("S12").Value = "  6    25"
